$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '66.464.94'
    'E2' = '  -1.86%  '
    'D3' = '3.838.63'
    'E3' = '  +1.70%  '
    'E4' = '  -0.38%  '
    'D5' = '421.82'
    'E5' = '  +0.29%  '
    'D6' = '127.95'
    'E6' = '  -3.68%  '
    'D7' = '3.836.67'
    'E7' = '  +2.04%  '
    'D8' = '0.602'
    'E8' = '  -7.58%  '
    'D9' = '0.998'
    'E9' = '  -0.10%  '
    'D10' = '0.716'
    'E10' = '  -7.63%  '
    'D11' = '0.163'
    'E11' = '  -14.42%  '
    'D12' = '0.0000348'
    'E12' = '  -18.88%  '
    'D13' = '40.10'
    'E13' = '  -6.58%  '
    'D14' = '4.429.32'
    'E14' = '  +1.22%  '
    'B15' = 'Polkadot'
    'C15' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D15' = '9.91'
    'E15' = '  -5.25%  '
    'B16' = 'Uniswap'
    'C16' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D16' = '15.81'
    'E16' = '  +19.86%  '
    'B17' = 'WrappedEther'
    'C17' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D17' = '3.847.32'
    'E17' = '  +2.11%  '
    'B18' = 'TRON'
    'C18' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D18' = '0.137'
    'E18' = '  -1.94%  '
    'D19' = '19.44'
    'E19' = '  -5.56%  '
    'D20' = '66.631.60'
    'E20' = '  -1.70%  '
    'E21' = '  -6.56%  '
    'D22' = '401.73'
    'E22' = '  -10.79%  '
    'D23' = '14.19'
    'E23' = '  -11.19%  '
    'D24' = '83.74'
    'E24' = '  -7.51%  '
    'D25' = '2.96'
    'E25' = '  -4.42%  '
    'D26' = '37.08'
    'E26' = '  -3.61%  '
    'D27' = '5.78'
    'E27' = '  +13.19%  '
    'D28' = '3.17'
    'E28' = '  -5.41%  '
    'D29' = '9.38'
    'E29' = '  -7.75%  '
    'D30' = '699.07'
    'E30' = '  +1.54%  '
    'B31' = 'Hedera'
    'C31' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D31' = '0.120'
    'E31' = '  -3.05%  '
    'B32' = 'Toncoin'
    'C32' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D32' = '2.75'
    'E32' = '  -0.53%  '
    'D33' = '12.23'
    'E33' = '  -3.96%  '
    'E34' = '  +3.13%  '
    'D35' = '0.149'
    'E35' = '  -10.46%  '
    'B36' = 'Dai'
    'C36' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D36' = '1.00'
    'E36' = '  +0.01%  '
    'B37' = 'InjectiveProtocol'
    'C37' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'D37' = '37.55'
    'E37' = '  -10.25%  '
    'D38' = '54.70'
    'E38' = '  -5.78%  '
    'D39' = '0.0₃0764'
    'E39' = '  +2.47%  '
    'D40' = '0.0451'
    'E40' = '  -8.59%  '
    'D41' = '2.91'
    'E41' = '  -3.91%  '
    'E42' = '  +0.12%  '
    'E43' = '  -9.72%  '
    'B44' = 'ApeXProtocol'
    'C44' = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
    'D44' = '3.15'
    'E44' = '  -1.27%  '
    'D45' = '3.30'
    'E45' = '  -3.04%  '
    'B46' = 'NEARProtocol'
    'C46' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D46' = '4.39'
    'E46' = '  +1.11%  '
    'B47' = 'Monero'
    'C47' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D47' = '143.50'
    'E47' = '  -3.21%  '
    'B48' = 'EnergySwap'
    'C48' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D48' = '26.29'
    'E48' = '  -8.18%  '
    'B49' = 'ARBITRUM'
    'C49' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D49' = '2.04'
    'E49' = '  -4.10%  '
    'E50' = '  -5.24%  '
    'D51' = '2.71'
    'E51' = '  -6.75%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}
